# Swap the Id/coordinates/time/substrate data between rows 10 and 11.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row1 = 10
$row2 = 11

# Columns whose values differ between the two rows and must be swapped.
$cols = @("A", "Q", "R", "Z", "AB", "AJ", "AK", "AO")

foreach ($col in $cols) {
    $rng1 = $ws.Range("$col$row1")
    $rng2 = $ws.Range("$col$row2")

    $val1 = $rng1.Value2
    $val2 = $rng2.Value2

    $rng1.Value2 = $val2
    $rng2.Value2 = $val1
}
